# Apply the TestData.xlsx edit:
#  - Replace the stray "C9 = 1" leftover row with a proper 4th data row
#    (LoginTest / firstName / Osanda) that matches the formatting used by
#    the other data rows.
#  - Update the current selection to N18 (cosmetic, matches saved file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (fill/border/font) of the last existing data row (row 3)
# onto the new row 4 so the new row matches the existing table styling.
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A4:C4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new data row.
$ws.Range("A4").Value = "LoginTest"
$ws.Range("B4").Value = "firstName"
$ws.Range("C4").Value = "Osanda"

# Drop the old leftover cell value that lived down at row 9.
$ws.Range("C9").ClearContents()

# Match the saved selection state.
$ws.Range("N18").Select() | Out-Null
